$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.240.90'
$ws.Range("E2").Value = '  +3.34%  '
$ws.Range("D3").Value = '1.738.96'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("D4").Value = '''0.9974'
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '''240.61'
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -1.75%  '
$ws.Range("D8").Value = '''0.2594'
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").Value = '''0.06151'
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '1.731.80'
$ws.Range("E10").Value = '  +1.40%  '
$ws.Range("D11").Value = '''16.08'
$ws.Range("E11").Value = '  +4.05%  '
$ws.Range("D12").Value = '''0.06935'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '''0.6025'
$ws.Range("E13").Value = '  +1.21%  '
$ws.Range("D14").Value = '''4.437'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '''76.90'
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("D16").Value = '''0.9986'
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = '27.192.93'
$ws.Range("E17").Value = '  +3.55%  '
$ws.Range("D18").Value = '''0.9973'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '''0.000007089'
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("E20").Value = '  +2.14%  '
$ws.Range("D21").Value = '1.951.13'
$ws.Range("E21").Value = '  +0.86%  '
$ws.Range("D22").Value = '''4.414'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("D23").Value = '''8.388'
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").Value = '''5.098'
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("D25").Value = '''142.09'
$ws.Range("E25").Value = '  +4.33%  '
$ws.Range("D26").Value = '''15.25'
$ws.Range("E26").Value = '  +0.75%  '
$ws.Range("D27").Value = '''1.817'
$ws.Range("E27").Value = '  +5.72%  '
$ws.Range("D28").Value = '''106.84'
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = '''1.384'
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("D30").Value = '''3.940'
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("D31").Value = '''0.07931'
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").Value = '''3.667'
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("D33").Value = '''0.04755'
$ws.Range("E33").Value = '  +7.15%  '
$ws.Range("D34").Value = '''2.600'
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("E35").Value = '  +2.23%  '
$ws.Range("D36").Value = '''0.6170'
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("D37").Value = '''0.9228'
$ws.Range("E37").Value = '  -2.35%  '
$ws.Range("D38").Value = '''2.537'
$ws.Range("E38").Value = '  +7.40%  '
$ws.Range("D39").Value = '''2.023'
$ws.Range("E39").Value = '  +2.23%  '
$ws.Range("D40").Value = '''0.9983'
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("D41").Value = '''5.693'
$ws.Range("E41").Value = '  +6.04%  '
$ws.Range("D42").Value = '''0.01487'
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("D43").Value = '''98.65'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").Value = '''0.3822'
$ws.Range("E44").Value = '  +1.08%  '
$ws.Range("D45").Value = '''6.845'
$ws.Range("E45").Value = '  +1.00%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").Value = '''0.05353'
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("D48").Value = '''7.793'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").Value = '''29.86'
$ws.Range("E49").Value = '  -1.56%  '
$ws.Range("D50").Value = '''1.242'
$ws.Range("E50").Value = '  +3.77%  '
$ws.Range("D51").Value = '''50.96'
$ws.Range("E51").Value = '  +0.20%  '
